$d = $word.ActiveDocument
$d.Content.Find.Execute("under 185 lines", $true, $false, $false, $false, $false,
                         $true, 1, $false, "around 230 lines", 2)
